$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M4").Value = 1115.12
$ws1.Range("M5").Value = 2690.98
$ws1.Range("I34").Value = 313.2
$ws1.Range("I55").Value = "5 de 53"
$ws1.Range("M55").Value = "7 de 53"

# --- Sheet 2: "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F4").Value = 3114.29
$ws2.Range("F5").Value = 4648.26
$ws2.Range("F34").Value = 4324.38
$ws2.Range("F59").Value = 17979.87

# --- Sheet 3: "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D7").Value = 877.5
$ws3.Range("E7").Value = 9.211016287574012
$ws3.Range("F7").Value = 0.9896121553489454

$ws3.Range("D12").Value = 8958.83
$ws3.Range("E12").Value = 52904.8903947566
$ws3.Range("F12").Value = 0.1448155711107107

$ws3.Range("D15").Value = 17459.46
$ws3.Range("E15").Value = 104595.3755108343
$ws3.Range("F15").Value = 0.1430460327681994
